$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.704964581735112
$ws.Cells.Item(2, 3).Value = 0.1941479203010203
$ws.Cells.Item(2, 4).Value = 0.07918298642189825
$ws.Cells.Item(2, 5).Value = 0.1122966986414156
$ws.Cells.Item(2, 7).Value = 0.002401136788963063
$ws.Cells.Item(2, 9).Value = 0.4389413944494009
$ws.Cells.Item(2, 13).Value = 0.3353924175294054
$ws.Cells.Item(2, 14).Value = 1.010228305635227
$ws.Cells.Item(2, 15).Value = 2.247061471130479

$ws.Cells.Item(3, 2).Value = 0.6211666978468315
$ws.Cells.Item(3, 3).Value = 0.1692864066073128
$ws.Cells.Item(3, 4).Value = 0.07173041045494699
$ws.Cells.Item(3, 5).Value = 0.1062102202841615
$ws.Cells.Item(3, 7).Value = 0.002404374914167823
$ws.Cells.Item(3, 9).Value = 0.4407376141199677
$ws.Cells.Item(3, 13).Value = 0.3002311413535779
$ws.Cells.Item(3, 14).Value = 1.02334869512595
$ws.Cells.Item(3, 15).Value = 2.219998014702497

$ws.Cells.Item(4, 2).Value = 0.5697016498697565
$ws.Cells.Item(4, 3).Value = 0.1539761553915753
$ws.Cells.Item(4, 4).Value = 0.06719039592093168
$ws.Cells.Item(4, 5).Value = 0.1025699261534143
$ws.Cells.Item(4, 7).Value = 0.0024064693446293
$ws.Cells.Item(4, 9).Value = 0.4421755267430427
$ws.Cells.Item(4, 13).Value = 0.2787351243650704
$ws.Cells.Item(4, 14).Value = 1.031864350281911
$ws.Cells.Item(4, 15).Value = 2.20514321083445

$ws.Cells.Item(5, 2).Value = 0.5487267501344206
$ws.Cells.Item(5, 3).Value = 0.1477258313397556
$ws.Cells.Item(5, 4).Value = 0.06534928775027993
$ws.Cells.Item(5, 5).Value = 0.1011105665215766
$ws.Cells.Item(5, 7).Value = 0.00240734963255456
$ws.Cells.Item(5, 9).Value = 0.4428455031918368
$ws.Cells.Item(5, 13).Value = 0.2699986101648477
$ws.Cells.Item(5, 14).Value = 1.035450067930157
$ws.Cells.Item(5, 15).Value = 2.199531523722698

$ws.Cells.Item(6, 2).Value = 0.545243756653889
$ws.Cells.Item(6, 3).Value = 0.1466872916415127
$ws.Cells.Item(6, 4).Value = 0.06504411457071058
$ws.Cells.Item(6, 5).Value = 0.1008696883907305
$ws.Cells.Item(6, 7).Value = 0.002407497424394341
$ws.Cells.Item(6, 9).Value = 0.4429618198789349
$ws.Cells.Item(6, 13).Value = 0.2685493212573604
$ws.Cells.Item(6, 14).Value = 1.036052448655127
$ws.Cells.Item(6, 15).Value = 2.198626348948949

$ws.Cells.Item(7, 2).Value = 0.5694187840637142
$ws.Cells.Item(7, 3).Value = 0.1538919068564439
$ws.Cells.Item(7, 4).Value = 0.06716552976813261
$ws.Cells.Item(7, 5).Value = 0.1025501475211712
$ws.Cells.Item(7, 7).Value = 0.002406481107882201
$ws.Cells.Item(7, 9).Value = 0.442184222411047
$ws.Cells.Item(7, 13).Value = 0.2786172065660182
$ws.Cells.Item(7, 14).Value = 1.031912240888779
$ws.Cells.Item(7, 15).Value = 2.205065742789174

$ws.Cells.Item(8, 2).Value = 0.6760741503084091
$ws.Cells.Item(8, 3).Value = 0.1855850933333727
$ws.Cells.Item(8, 4).Value = 0.07660586011290604
$ws.Cells.Item(8, 5).Value = 0.1101778264492523
$ws.Cells.Item(8, 7).Value = 0.002402231301475204
$ws.Cells.Item(8, 9).Value = 0.439491057000911
$ws.Cells.Item(8, 13).Value = 0.3232493832441747
$ws.Cells.Item(8, 14).Value = 1.014656765584554
$ws.Cells.Item(8, 15).Value = 2.237363391814739

$ws.Cells.Item(9, 2).Value = 0.8851010841498237
$ws.Cells.Item(9, 3).Value = 0.2473759845999552
$ws.Cells.Item(9, 4).Value = 0.095406028847421
$ws.Cells.Item(9, 5).Value = 0.1259164861639803
$ws.Cells.Item(9, 7).Value = 0.002394736343362596
$ws.Cells.Item(9, 9).Value = 0.4368788710264297
$ws.Cells.Item(9, 13).Value = 0.4115228065738421
$ws.Cells.Item(9, 14).Value = 0.9844691588410797
$ws.Cells.Item(9, 15).Value = 2.314750478456688

$ws.Cells.Item(10, 2).Value = 1.038583172769052
$ws.Cells.Item(10, 3).Value = 0.2925580362645235
$ws.Cells.Item(10, 4).Value = 0.1093993153735511
$ws.Cells.Item(10, 5).Value = 0.1379742414340939
$ws.Cells.Item(10, 7).Value = 0.002389735864297102
$ws.Cells.Item(10, 9).Value = 0.436602485730802
$ws.Cells.Item(10, 13).Value = 0.4768560257221139
$ws.Cells.Item(10, 14).Value = 0.9645193171226012
$ws.Cells.Item(10, 15).Value = 2.380277928565789

$ws.Cells.Item(11, 2).Value = 1.108384562073127
$ws.Cells.Item(11, 3).Value = 0.3130666915055258
$ws.Cells.Item(11, 4).Value = 0.1158057277175146
$ws.Cells.Item(11, 5).Value = 0.1435711011478773
$ws.Cells.Item(11, 7).Value = 0.002387569765615347
$ws.Cells.Item(11, 9).Value = 0.4368368856457678
$ws.Cells.Item(11, 13).Value = 0.5066868274648328
$ws.Cells.Item(11, 14).Value = 0.9559284844248133
$ws.Cells.Item(11, 15).Value = 2.411994188282677

$ws.Cells.Item(12, 2).Value = 1.134813368723485
$ws.Cells.Item(12, 3).Value = 0.3208263126389568
$ws.Cells.Item(12, 4).Value = 0.1182376038815107
$ws.Cells.Item(12, 5).Value = 0.145706858625239
$ws.Cells.Item(12, 7).Value = 0.002386765057018136
$ws.Cells.Item(12, 9).Value = 0.4369777010064055
$ws.Cells.Item(12, 13).Value = 0.5179991149250327
$ws.Cells.Item(12, 14).Value = 0.9527451265222204
$ws.Cells.Item(12, 15).Value = 2.424280296073647

$ws.Cells.Item(13, 2).Value = 1.129121619448995
$ws.Cells.Item(13, 3).Value = 0.3191554329013115
$ws.Cells.Item(13, 4).Value = 0.1177135924881298
$ws.Cells.Item(13, 5).Value = 0.1452461534420095
$ws.Cells.Item(13, 7).Value = 0.002386937675205578
$ws.Cells.Item(13, 9).Value = 0.4369450545575617
$ws.Cells.Item(13, 13).Value = 0.515562097167404
$ws.Cells.Item(13, 14).Value = 0.9534276127584995
$ws.Cells.Item(13, 15).Value = 2.421621969709747

$ws.Cells.Item(14, 2).Value = 1.110558950023858
$ws.Cells.Item(14, 3).Value = 0.3137052123162789
$ws.Cells.Item(14, 4).Value = 0.1160056810205816
$ws.Cells.Item(14, 5).Value = 0.1437464819380097
$ws.Cells.Item(14, 7).Value = 0.002387503250652037
$ws.Cells.Item(14, 9).Value = 0.4368474258542179
$ws.Cells.Item(14, 13).Value = 0.5076171745057678
$ws.Cells.Item(14, 14).Value = 0.9556651878550575
$ws.Cells.Item(14, 15).Value = 2.412999436436053

$ws.Cells.Item(15, 2).Value = 1.099188303536721
$ws.Cells.Item(15, 3).Value = 0.3103659371857077
$ws.Cells.Item(15, 4).Value = 0.1149603063753943
$ws.Cells.Item(15, 5).Value = 0.1428300272517831
$ws.Cells.Item(15, 7).Value = 0.002387851704285738
$ws.Cells.Item(15, 9).Value = 0.4367944122567096
$ws.Cells.Item(15, 13).Value = 0.502752769053842
$ws.Cells.Item(15, 14).Value = 0.9570448620101217
$ws.Cells.Item(15, 15).Value = 2.407753855637196

$ws.Cells.Item(16, 2).Value = 1.03402104760039
$ws.Cells.Item(16, 3).Value = 0.2912168350113689
$ws.Cells.Item(16, 4).Value = 0.1089814656396442
$ws.Cells.Item(16, 5).Value = 0.1376107457855511
$ws.Cells.Item(16, 7).Value = 0.002389879603433331
$ws.Cells.Item(16, 9).Value = 0.436594436675982
$ws.Cells.Item(16, 13).Value = 0.4749087456198566
$ws.Cells.Item(16, 14).Value = 0.9650905085280002
$ws.Cells.Item(16, 15).Value = 2.378243712732427

$ws.Cells.Item(17, 2).Value = 0.99403775272026
$ws.Cells.Item(17, 3).Value = 0.2794579080366759
$ws.Cells.Item(17, 4).Value = 0.1053241181624287
$ws.Cells.Item(17, 5).Value = 0.1344377016679346
$ws.Cells.Item(17, 7).Value = 0.00239115142548743
$ws.Cells.Item(17, 9).Value = 0.4365641856230411
$ws.Cells.Item(17, 13).Value = 0.4578557092745683
$ws.Cells.Item(17, 14).Value = 0.9701504298832937
$ws.Cells.Item(17, 15).Value = 2.360629854635818

$ws.Cells.Item(18, 2).Value = 0.9710387375336609
$ws.Cells.Item(18, 3).Value = 0.2726902656120842
$ws.Cells.Item(18, 4).Value = 0.1032243461457796
$ws.Cells.Item(18, 5).Value = 0.1326231498790875
$ws.Cells.Item(18, 7).Value = 0.002391893174478546
$ws.Cells.Item(18, 9).Value = 0.436580667041504
$ws.Cells.Item(18, 13).Value = 0.448057636221435
$ws.Cells.Item(18, 14).Value = 0.9731063528331063
$ws.Cells.Item(18, 15).Value = 2.350678278525635

$ws.Cells.Item(19, 2).Value = 0.9632514024764305
$ws.Cells.Item(19, 3).Value = 0.270398136893391
$ws.Cells.Item(19, 4).Value = 0.102514057283841
$ws.Cells.Item(19, 5).Value = 0.1320105685097772
$ws.Cells.Item(19, 7).Value = 0.002392146077444887
$ws.Cells.Item(19, 9).Value = 0.436592058341688
$ws.Cells.Item(19, 13).Value = 0.4447419561244459
$ws.Cells.Item(19, 14).Value = 0.9741150034874408
$ws.Cells.Item(19, 15).Value = 2.347339619178683

$ws.Cells.Item(20, 2).Value = 0.9982942227112517
$ws.Cells.Item(20, 3).Value = 0.2807101032807964
$ws.Cells.Item(20, 4).Value = 0.1057130515838196
$ws.Cells.Item(20, 5).Value = 0.1347743889597837
$ws.Cells.Item(20, 7).Value = 0.002391014979593951
$ws.Cells.Item(20, 9).Value = 0.4365638974055841
$ws.Cells.Item(20, 13).Value = 0.4596699585905668
$ws.Cells.Item(20, 14).Value = 0.9696070732141351
$ws.Cells.Item(20, 15).Value = 2.362486296334197

$ws.Cells.Item(21, 2).Value = 1.116011356546949
$ws.Cells.Item(21, 3).Value = 0.315306253139596
$ws.Cells.Item(21, 4).Value = 0.1165071753870706
$ws.Cells.Item(21, 5).Value = 0.1441865258711914
$ws.Cells.Item(21, 7).Value = 0.002387336706216521
$ws.Cells.Item(21, 9).Value = 0.436874686929734
$ws.Cells.Item(21, 13).Value = 0.5099503552705045
$ws.Cells.Item(21, 14).Value = 0.9550060622031751
$ws.Cells.Item(21, 15).Value = 2.415524584618595

$ws.Cells.Item(22, 2).Value = 1.192925736743064
$ws.Cells.Item(22, 3).Value = 0.3378785429185598
$ws.Cells.Item(22, 4).Value = 0.1235962075739963
$ws.Cells.Item(22, 5).Value = 0.1504333230971255
$ws.Cells.Item(22, 7).Value = 0.002385023319791582
$ws.Cells.Item(22, 9).Value = 0.4373813362068759
$ws.Cells.Item(22, 13).Value = 0.5429048859361529
$ws.Cells.Item(22, 14).Value = 0.9458703973839562
$ws.Cells.Item(22, 15).Value = 2.451796752753694

$ws.Cells.Item(23, 2).Value = 1.151877236831979
$ws.Cells.Item(23, 3).Value = 0.3258348259033994
$ws.Cells.Item(23, 4).Value = 0.1198094923338431
$ws.Cells.Item(23, 5).Value = 0.147090465165121
$ws.Cells.Item(23, 7).Value = 0.002386249755022495
$ws.Cells.Item(23, 9).Value = 0.4370830664684462
$ws.Cells.Item(23, 13).Value = 0.5253078424501325
$ws.Cells.Item(23, 14).Value = 0.9507089908364996
$ws.Cells.Item(23, 15).Value = 2.432289914336792

$ws.Cells.Item(24, 2).Value = 0.9963699102164583
$ws.Cells.Item(24, 3).Value = 0.2801440084795388
$ws.Cells.Item(24, 4).Value = 0.1055372058008004
$ws.Cells.Item(24, 5).Value = 0.1346221425269292
$ws.Cells.Item(24, 7).Value = 0.002391076633866604
$ws.Cells.Item(24, 9).Value = 0.4365639222125068
$ws.Cells.Item(24, 13).Value = 0.4588497179744877
$ws.Cells.Item(24, 14).Value = 0.9698525786551713
$ws.Cells.Item(24, 15).Value = 2.361646454451289

$ws.Cells.Item(25, 2).Value = 0.8285685403780576
$ws.Cells.Item(25, 3).Value = 0.2306980408184245
$ws.Cells.Item(25, 4).Value = 0.09028872449209757
$ws.Cells.Item(25, 5).Value = 0.1215732577294162
$ws.Cells.Item(25, 7).Value = 0.002396674672865943
$ws.Cells.Item(25, 9).Value = 0.4372981778643066
$ws.Cells.Item(25, 13).Value = 0.3875600049020278
$ws.Cells.Item(25, 14).Value = 0.9844691588410797
$ws.Cells.Item(25, 15).Value = 2.314750478456688
